# Simulated Wild Card round and logged it.
# Updates the Rushing and Receiving stat sheets with the new cumulative
# totals produced by the Wild Card game.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet -------------------------------------------------------
$wsRush = $wb.Worksheets.Item("Rushing")

# J.Mixon (row 4)
$wsRush.Range("C4").Value = 192
$wsRush.Range("D4").Value = 95
$wsRush.Range("E4").Value = 20
$wsRush.Range("F4").Value = 41

# S.Perine (row 5)
$wsRush.Range("C5").Value = 24

# C.Evans (row 6)
$wsRush.Range("D6").Value = 4

# T.Boyd (row 8)
$wsRush.Range("D8").Value = 1

# J.Chase (row 9)
$wsRush.Range("D9").Value = 3
$wsRush.Range("E9").Value = 2

# --- Receiving sheet ------------------------------------------------------
$wsRecv = $wb.Worksheets.Item("Receiving")

# J.Mixon (row 2)
$wsRecv.Range("C2").Value = 49
$wsRecv.Range("D2").Value = 42
$wsRecv.Range("E2").Value = 4
$wsRecv.Range("F2").Value = 4
$wsRecv.Range("G2").Value = 8
$wsRecv.Range("H2").Value = 6

# T.Boyd (row 5)
$wsRecv.Range("C5").Value = 81
$wsRecv.Range("D5").Value = 61
$wsRecv.Range("G5").Value = 8
$wsRecv.Range("H5").Value = 6

# T.Higgins (row 6)
$wsRecv.Range("C6").Value = 92
$wsRecv.Range("D6").Value = 66
$wsRecv.Range("E6").Value = 33
$wsRecv.Range("G6").Value = 11

# J.Chase (row 7)
$wsRecv.Range("C7").Value = 94
$wsRecv.Range("D7").Value = 67
$wsRecv.Range("E7").Value = 43
$wsRecv.Range("F7").Value = 22

# C.Uzomah (row 11)
$wsRecv.Range("C11").Value = 60
$wsRecv.Range("D11").Value = 47
$wsRecv.Range("E11").Value = 9
$wsRecv.Range("F11").Value = 8
$wsRecv.Range("G11").Value = 4
$wsRecv.Range("H11").Value = 2
